# Data source corrected and updated
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns J and K held stray/incorrect data: row 1 had the leftover text
# labels "r"/"s" while rows 2-51 held 0.3 / 0.6. The corrected data source
# uses a constant 0.6 / 1 down the full J1:K51 range.
$ws.Range("J1:J51").Value = 0.6
$ws.Range("K1:K51").Value = 1

# Reflect the new viewport/selection saved with the corrected sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 41
$win.ScrollColumn = 2
$ws.Range("K1:K51").Select()
